# Natmi following Dr Hou advice
# Update LR-pair rows: add 'ECs' target cluster, refresh stats, extend to 6 data rows (2-7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ A="FAPs"; B="Wnt5a"; C="Fzd2"; D="ECs"; N=@(3,1,9.156959333333335,27.470878,0.969469463764299,0.9694694637642989,1,0.3333333333333333,0.1278803333333333,0.383641,0.009974564977605908,0.009974564977605908,1.170995011866445,10.538955106798,0.009670036160121756,0.009670036160121754) },
    @{ A="FAPs"; B="Wnt5a"; C="Fzd2"; D="FAPs"; N=@(3,1,9.156959333333335,27.470878,0.969469463764299,0.9694694637642989,3,1,10.61985133333333,31.859554,0.8283400145723324,0.8283400145723324,97.24554678537912,875.2099210684121,0.8030503497419507,0.8030503497419507) },
    @{ A="FAPs"; B="Wnt5a"; C="Fzd2"; D="sCs"; N=@(3,1,9.156959333333335,27.470878,0.969469463764299,0.9694694637642989,3,1,2.072911,6.218733,0.1616854204500617,0.1616854204500617,18.98156172861934,170.834055557574,0.1567490778622265,0.1567490778622265) },
    @{ A="sCs"; B="Wnt5a"; C="Fzd2"; D="ECs"; N=@(3,1,0.288371,0.865113,0.03053053623570109,0.03053053623570109,1,0.3333333333333333,0.1278803333333333,0.383641,0.009974564977605908,0.009974564977605908,0.03687697960366667,0.331892816433,0.0003045288174841522,0.0003045288174841522) },
    @{ A="sCs"; B="Wnt5a"; C="Fzd2"; D="FAPs"; N=@(3,1,0.288371,0.865113,0.03053053623570109,0.03053053623570109,3,1,10.61985133333333,31.859554,0.8283400145723324,0.8283400145723324,3.062457148844667,27.562114339602,0.02528966483038176,0.02528966483038176) },
    @{ A="sCs"; B="Wnt5a"; C="Fzd2"; D="sCs"; N=@(3,1,0.288371,0.865113,0.03053053623570109,0.03053053623570109,3,1,2.072911,6.218733,0.1616854204500617,0.1616854204500617,0.597767417981,5.379906761829,0.004936342587835174,0.004936342587835174) }
)

$r = 2
foreach ($row in $rowData) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    for ($i = 0; $i -lt $row.N.Length; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $row.N[$i]
    }
    $r++
}
